$wb = $excel.ActiveWorkbook
$wsNos = $wb.Worksheets.Item("Nós")
$wsArestas = $wb.Worksheets.Item("Arestas")

# --- Sheet "Nos": append new scientist names (rows 29-49) ---
$sheet1NewNames = @(
    "Isaac Newton",
    "John Freind",
    "John Kiel",
    "Robert Boyle",
    "Daniel Bernoulli",
    "John Herapath",
    "John James Waterson",
    "Antoine Lavoisier ",
    "John Dalton",
    "Joseph Priestley",
    "Charles Blagden",
    "Henry Cavendish",
    "William Henry",
    "Joseph Black",
    "André-Marie Ampère ",
    "Michael Faraday",
    "Heirich Hertz",
    "Amedeo Avogadro ",
    "Jacob Berzellius",
    "Henry Becquerel",
    "James Chadwick",
)

$startRow = 29
for ($i = 0; $i -lt $sheet1NewNames.Count; $i++) {
    $r = $startRow + $i
    $wsNos.Cells.Item($r, 1).Value = $sheet1NewNames[$i]
}
# copy row-28 formatting down onto the new rows (reuses the existing shared style)
$wsNos.Range("A28").Copy()
$wsNos.Range("A29:A49").PasteSpecial(-4122)

# --- Sheet "Arestas": rename headers, add Seculo_Interacao column, append new edges ---
$wsArestas.Range("A1").Value = "Cientista_Origem"
$wsArestas.Range("B1").Value = "Cientista_Destino"
$wsArestas.Range("C1").Value = "Seculo_Interacao"

# style the brand-new C1 header cell off of the existing A1/B1 style, then restore its text
$wsArestas.Range("A1").Copy()
$wsArestas.Range("C1").PasteSpecial(-4122)
$wsArestas.Range("C1").Value = "Seculo_Interacao"

# set the new column width, then fill the century value for the pre-existing rows (2-28)
$wsArestas.Columns.Item(3).ColumnWidth = 14.5
for ($r = 2; $r -le 28; $r++) {
    $wsArestas.Cells.Item($r, 3).Value = 20
}

# append the new edge rows (30 historical interactions, 19th/18th century)
$sheet2NewRows = @(
    @("John Dalton", "William Henry", 19),
    @("John Joseph Thomsom", "William Crookes", 19),
    @("John Joseph Thomsom", "Heinrich Hertz", 19),
    @("Michael Faraday", "Humphry Davy", 19),
    @("Humphry Davy", "William Wollaston", 19),
    @("Louis J. Gay-Lussac", "John Dalton", 19),
    @("Louis J. Gay-Lussac", "Amedeo Avogadro", 19),
    @("Claude Louis Berthollet", "Louis J. Gay-Lussac", 19),
    @("Humphry Davy", "Jacob Berzellius", 19),
    @("Jacob Berzellius ", "John Dalton", 19),
    @("Robert Brown", "Giovanni Cantoni", 19),
    @("Robert Brown ", "Bodoszewski", 19),
    @("Max Planck ", "Albert Einstein", 19),
    @("Albert Einstein", "Henry Becquerel", 19),
    @("Ernest Rutherford", "Marie Curie", 19),
    @("Ernest Rutherford", "James Chadwick", 19),
    @("Isaac Newton", "John Freind", 18),
    @("Isaac Newton", "John Kiel", 18),
    @("Robert Boyle", "Isaac Newton", 18),
    @("Daniel Bernoulli ", "John Herapath", 18),
    @("Daniel Bernoulli ", "John James Waterson", 18),
    @("Antoine Lavoisier", "John Dalton", 18),
    @("Antoine Lavoisier ", "Joseph Priestley", 18),
    @("Charles Blagden", "Henry Cavendish", 18),
    @("Henry Cavendish", "Joseph Priestley", 18),
    @("Joseph Priestley", "William Henry", 18),
    @("Joseph Priestley", "Joseph Black", 18),
    @("André-Marie Àmpere", "Robert Boyle", 18),
    @("Joseph Priestley", "John Dalton", 18),
    @("Robert Boyle", "John Dalton", 18),
)

$startRow2 = 29
for ($i = 0; $i -lt $sheet2NewRows.Count; $i++) {
    $r = $startRow2 + $i
    $row = $sheet2NewRows[$i]
    $wsArestas.Cells.Item($r, 1).Value = $row[0]
    $wsArestas.Cells.Item($r, 2).Value = $row[1]
    $wsArestas.Cells.Item($r, 3).Value = $row[2]
}

# copy the (now fully styled) header row formatting down across the whole used range
# so every new A/B/C data cell shares the same style as the rest of the table
$wsArestas.Range("A1:C1").Copy()
$wsArestas.Range("A2:C58").PasteSpecial(-4122)

$wsArestas.Range("A1").Select()
